$wb = $excel.ActiveWorkbook

# --- Update header text on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match page margins used on the other sheets (0.75"/0.75"/1"/1"/0.5"/0.5")
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Header row text
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$data = @(
    @(44934.99999999999, 30, 30.00008492383639, 30.00008492394558),
    @(44941.99999999999, 10, 10.00007495009424, 10.00007495020469),
    @(44948.99999999999, 0, -9.999935269997025, -9.999934796718197),
    @(44955.99999999999, 0, -29.99994577640831, -29.9999442358682),
    @(44962.99999999999, 0, -49.99995641903415, -49.99995354564264),
    @(44969.99999999999, 0, -69.99996726973102, -69.99996277872472),
    @(44976.99999999999, 0, -89.99997817756405, -89.99997185193499),
    @(44983.99999999999, 0, -109.9999891685134, -109.9999806718697),
    @(44990.99999999999, 0, -130.0000003741229, -129.9999893895295),
    @(44997.99999999999, 0, -150.0000115606938, -149.9999980046995)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# --- Match formatting used on the other sheets: bold/bordered/centered
# header style, and the date number format on column A ---
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A11").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Restore original active sheet selection
$wsWeekly.Activate()
